# Fill in the "Absent" (column H) values that were left blank/incorrect
# when the consolidated report was formed. H = 1 - D (Absent = NOT Present),
# matching the logic already used for the already-populated H cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H14").Value = 0
